$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header text updates (component renamed to be more specific) ---
$ws.Range("A7").Value = "MCP3008 ADC"
$ws.Range("A12").Value = "DS3231 RTC"
$ws.Range("A15").Value = "ATTiny84 Microcontroller"
$ws.Range("A23").Value = "L297 Motor Control"

# --- Insert a new row for the "Limit Switches" section header, pushing the ---
# --- two existing limit-switch pin rows down from 29/30 to 30/31          ---
$ws.Rows(29).Insert()

# --- Fill in the motor-control (L297) rows that Mike had left blank ---
$ws.Range("A24").Value = "L297_1"
$ws.Range("B24").Value = "AZIMUTH_EN"
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = " enable for azimuth motor driver"

$ws.Range("A25").Value = "L297_2"
$ws.Range("B25").Value = "ELEVATION_EN"
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 27
$ws.Range("E25").Value = " enable for elevation motor driver "

$ws.Range("A26").Value = "L297"
$ws.Range("B26").Value = "RESET_N"
$ws.Range("C26").Value = 31
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = " shared between both L297 driver chips "

$ws.Range("A27").Value = "L297"
$ws.Range("B27").Value = "CLOCK_N"
$ws.Range("C27").Value = 21
$ws.Range("D27").Value = 9
$ws.Range("E27").Value = " shared between both L297 driver chips"

$ws.Range("A28").Value = "L297"
$ws.Range("B28").Value = "DIRECTION"
$ws.Range("C28").Value = 28
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = " shared between both L297 driver chips"

# --- New "Limit Switches" section header row, formatted like the other ---
# --- mid-sheet section banners (copy formatting from the Motor Control banner) ---
$ws.Range("A23:E23").Copy()
$ws.Range("A29:E29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("A29").Value = "Limit Switches"

# --- Label the two limit switch pin rows (values already shifted down by the insert) ---
$ws.Range("A30").Value = "LIM_SW"
$ws.Range("B30").Value = 0
$ws.Range("A31").Value = "LIM_SW"
$ws.Range("B31").Value = 0
